$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right after "总计" (i.e. before
#    what is currently the second sheet, "2022-Q3"). All the quarter sheets
#    after it shift one slot to the right; their content is left untouched.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3sheet    = $wb.Worksheets.Item(2)
$newSheet   = $wb.Worksheets.Add($q3sheet, $null)
$newSheet.Name = "2022-Q4"

# Header row (matches the other quarter sheets).
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160

# Single data row: 013623 / 湘财周期轮动一年持有期混合.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("A2").HorizontalAlignment = -4108
$newSheet.Range("A2").VerticalAlignment = -4160

$newSheet.Range("B2").Value = "'013623"
$newSheet.Range("C2").Value = "湘财周期轮动一年持有期混合"
$newSheet.Range("D2").Value = "'4.11"
$newSheet.Range("E2").Value = "'79.16"
$newSheet.Range("F2").Value = "'3.85"
$newSheet.Range("G2").Value = "'0.1582"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------------
# 2) Add a "2022-Q4" row at the top of the summary sheet "总计" (new row 2),
#    pushing the existing quarters down by one row and renumbering the index
#    column (A) sequentially.
# ---------------------------------------------------------------------------
$ws = $totalSheet
$ws.Rows("2:2").Insert()

# Copy A3's style (bold/centered index style) onto the freshly inserted A2,
# then clear the border/format that Insert() copied onto B2:D2 from the
# header row so the new row matches the look of the other data rows.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.16

# Renumber the index column and re-assert clean numeric literals for the
# shifted rows (avoids float noise introduced by the row-insert shift).
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2022-Q3"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0.02

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2021-Q3"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0.36

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2021-Q2"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 1.37

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "2021-Q1"
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 1.83

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "2020-Q4"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0.01

# Selection back on the summary sheet, matching the original workbook state.
$ws.Range("A1").Select()
